$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 534, carrying a copy of row 534's current
# contents/formatting into both the new row 534 and the (now) row 535 --
# this mirrors Excel's "Copy row, Insert Copied Cells" gesture.
$ws.Rows.Item(534).Copy()
$ws.Rows.Item(534).Insert()

# Now edit the freshly inserted row 534 with its own data: a new date
# (Fecha) and a new Origen.
$ws.Cells.Item(534, 4).Value = 45258
$ws.Cells.Item(534, 15).Value = "Provincia del Elquí"
